$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" everywhere it appears ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws2.Range("C2").Value = "In Translation"
$ws3.Range("C2").Value = "In Translation"

# --- Shrink the "Status" columns to match the regenerated report layout ---
# Target stored width is 13.4101848602295 characters; the COM ColumnWidth setter
# only supports discrete increments of 1/6 character, so we choose the closest
# representable value (12.5 -> stored width 13.3333...).
$newWidth = 12.5
$ws1.Columns.Item(5).ColumnWidth = $newWidth   # Overview!E (zh-cn status)
$ws1.Columns.Item(6).ColumnWidth = $newWidth   # Overview!F (de-de status)
$ws2.Columns.Item(3).ColumnWidth = $newWidth   # zh-cn!C (Status)
$ws3.Columns.Item(3).ColumnWidth = $newWidth   # de-de!C (Status)
